$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert 9 new rows at the top of the data (above the previous A2) to make room
# for the 9 newest epochs (501 down to 493), pushing the rest of the log down.
$ws.Range("A2:A10").Insert()

# Fill the newly inserted rows with the latest epoch results
$ws.Range("A2").Value = "Epoch:501, time:9.649806, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
$ws.Range("A3").Value = "Epoch:500, time:9.438108, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
$ws.Range("A4").Value = "Epoch:499, time:9.327127, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
$ws.Range("A5").Value = "Epoch:498, time:9.611853, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
$ws.Range("A6").Value = "Epoch:497, time:9.535737, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
$ws.Range("A7").Value = "Epoch:496, time:9.260726, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
$ws.Range("A8").Value = "Epoch:495, time:14.170702, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
$ws.Range("A9").Value = "Epoch:494, time:9.904561, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
$ws.Range("A10").Value = "Epoch:493, time:9.826703, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
